$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 371
$ws.Range("I2").Value = 1062
$ws.Range("J2").Value = 4455
$ws.Range("K2").Value = 23
$ws.Range("L2").Value = 1239
$ws.Range("M2").Value = 66
$ws.Range("N2").Value = 799
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 19
$ws.Range("Q2").Value = 8
$ws.Range("R2").Value = 55
$ws.Range("S2").Value = 458
$ws.Range("T2").Value = 756
$ws.Range("U2").Value = 58
$ws.Range("V2").Value = 7006
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 6876
$ws.Range("Z2").Value = 97
